# Update template with new infrastructure
# - Add a new "Orthography" column (L) with a value for the first data row
# - Widen the existing columns to fit the new, longer content
# - Move the active selection to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header and data
$ws.Range("L1").Value = "Orthography"
$ws.Range("L2").Value = "p/general"

# Widen columns (A..L and beyond) to match the new, wider template layout
$ws.Columns.Item(1).ColumnWidth = 25.5561224489796
$ws.Columns.Item(2).ColumnWidth = 25.5561224489796
$ws.Columns.Item(3).ColumnWidth = 28.2602040816327
$ws.Columns.Item(4).ColumnWidth = 14.7602040816327
$ws.Columns.Item(5).ColumnWidth = 17.6377551020408
$ws.Columns.Item(6).ColumnWidth = 98.4540816326531
$ws.Columns.Item(7).ColumnWidth = 18.719387755102
$ws.Columns.Item(8).ColumnWidth = 7.19897959183674
$ws.Columns.Item(9).ColumnWidth = 13.5
$ws.Columns.Item(10).ColumnWidth = 18.719387755102
$ws.Columns.Item(11).ColumnWidth = 16.7397959183673
$ws.Columns.Item(12).ColumnWidth = 20.6989795918367

# Move the selection to A2, as in the updated template
$ws.Range("A2").Select()
